# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 5.5

# Row 4
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5

# Row 7
$ws.Range("G7").Value = 4.33
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = 4.75
$ws.Range("L7").Value = 2.75
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("W7").Value = 9
$ws.Range("X7").Value = 19
$ws.Range("AF7").Value = 81
$ws.Range("AG7").Value = 5.5
$ws.Range("AH7").Value = 8
$ws.Range("AJ7").Value = 17
$ws.Range("AP7").Value = 41
$ws.Range("AU7").Value = 9.5
$ws.Range("AV7").Value = 81
$ws.Range("AW7").Value = 3.75
$ws.Range("AY7").Value = 26
$ws.Range("BA7").Value = 67

$wb.Save()
